$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (shifts existing rows 28-46 down to 29-47,
# matching the data shift described by the diff) and populate it with the
# new weekly price observation.
$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44488
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112052
$ws.Cells.Item(28, 7).Value = "Albahaca"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 800
$ws.Cells.Item(28, 11).Value = 3500
$ws.Cells.Item(28, 12).Value = 4000
$ws.Cells.Item(28, 13).Value = 3750
$ws.Cells.Item(28, 14).Value = "`$/paquete"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 3750
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = "Hortaliza"
